# Benchmark.xlsx update - "Add files via upload"
# Adds 10 new benchmark rows (QueryNum 21-30) to Sheet1, switches the
# workbook's base font from Arial to Calibri, and refreshes the window
# selection / scroll position to reflect the newly-added data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Re-font the existing data (Arial -> Calibri) before appending the
#    new rows, so the new rows naturally pick up the Calibri font too.
# ---------------------------------------------------------------------
$ws.Range("A1:E21").Font.Name = "Calibri"

# ---------------------------------------------------------------------
# 2) Append the new query rows (22-31 / QueryNum 21-30).
#    Column A = QueryNum, B = Query text, C = project, D = path (blank),
#    E = blank placeholder cell (present only on rows that keep the
#    trailing styled empty cell).
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row = 22; Num = 21; Query = "list iterable node"; Project = "IterableList"; HasE = $true },
    @{ Row = 23; Num = 22; Query = "long list with iteration option and the list must consisting of nodes"; Project = "IterableList"; HasE = $true },
    @{ Row = 24; Num = 23; Query = "execute xor operation between two blocks"; Project = "AES attack"; HasE = $true },
    @{ Row = 25; Num = 24; Query = "encryption and dectyption by AES attack"; Project = "AES attack"; HasE = $true },
    @{ Row = 26; Num = 25; Query = "solve maze game by bread search"; Project = "maze"; HasE = $true },
    @{ Row = 27; Num = 26; Query = "generate a maze that contains stat position and search function"; Project = "maze"; HasE = $true },
    @{ Row = 28; Num = 27; Query = "get value of binary expression"; Project = "Binary-Operations"; HasE = $true },
    @{ Row = 29; Num = 28; Query = "show solution of maze game by showing the path"; Project = "maze"; HasE = $true },
    @{ Row = 30; Num = 29; Query = "write descryption to a file"; Project = "AES attack"; HasE = $false },
    @{ Row = 31; Num = 30; Query = "write descryption to a file"; Project = "AES attack"; HasE = $false }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    $a = $ws.Cells.Item($rowNum, 1)
    $a.Value = $r.Num

    $b = $ws.Cells.Item($rowNum, 2)
    $b.Value = $r.Query
    $b.HorizontalAlignment = -4131   # xlLeft
    $b.VerticalAlignment = -4160     # xlTop
    $b.Font.Name = "Calibri"

    $c = $ws.Cells.Item($rowNum, 3)
    $c.Value = $r.Project
    $c.VerticalAlignment = -4160     # xlTop
    $c.Font.Name = "Calibri"

    if ($r.HasE) {
        $e = $ws.Cells.Item($rowNum, 5)
        $e.HorizontalAlignment = -4108   # xlCenter
        $e.VerticalAlignment = -4160     # xlTop
        $e.Font.Name = "Calibri"
    }
}

# ---------------------------------------------------------------------
# 3) Refresh the view: scroll so row 10 is at the top and select C32
#    (mirrors the author re-saving after scrolling through the sheet).
# ---------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("C32").Select()
